# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
#
# This script appends the newly logged Week 17 per-play / per-game figures
# to the running season log strings (YDS and ST sheets), and bumps the
# season-to-date cumulative totals (OFF, DEF, ST, TURNS, PEN sheets) by the
# Week 17 deltas.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet - append Week 17 per-play yardage logs
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = $ws.Range("B2").Value() + " 6 2 4 1 1 2 3 14 1 7 7 -2 11 4 7 1 8 1 -4"
$ws.Range("B3").Value = $ws.Range("B3").Value() + " 10 18 7 19 6 4 9 22 11 17 11 18 2 35 12 9 2 8 3 21 15 24 4 15 7"
$ws.Range("C2").Value = $ws.Range("C2").Value() + " 2 1 2 -3 5 3 7 1 1 1 9 5 17 11 3 -1 12 2 4 3 3 32 6 6 4 6 3 2 3 -2"
$ws.Range("C3").Value = $ws.Range("C3").Value() + " 15 13 5 1 14 9 18 5 33 2 14 15 4 15 7 11 4 5"

# ---------------------------------------------------------------------
# ST sheet - append Week 17 special-teams per-game logs, bump totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")

$ws.Range("B4").Value = $ws.Range("B4").Value() + " 66"
$ws.Range("B5").Value = $ws.Range("B5").Value() + " 22"
$ws.Range("B6").Value = $ws.Range("B6").Value() + " 18"
$ws.Range("D3").Value = $ws.Range("D3").Value() + " 39 53"
$ws.Range("D4").Value = $ws.Range("D4").Value() + " 0 0"
$ws.Range("D5").Value = $ws.Range("D5").Value() + " 0 0 18"

$ws.Range("B2").Value = $ws.Range("B2").Value() + 4
$ws.Range("D2").Value = $ws.Range("D2").Value() + 2
$ws.Range("F2").Value = $ws.Range("F2").Value() + 2
$ws.Range("G2").Value = $ws.Range("G2").Value() + 2
$ws.Range("H2").Value = $ws.Range("H2").Value() + 1
$ws.Range("N2").Value = $ws.Range("N2").Value() + 1
$ws.Range("B3").Value = $ws.Range("B3").Value() + 3

# ---------------------------------------------------------------------
# OFF sheet - bump Home (row 2) and Road (row 3) season totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")

$ws.Range("C2").Value = $ws.Range("C2").Value() + 10
$ws.Range("F2").Value = $ws.Range("F2").Value() + 1
$ws.Range("G2").Value = $ws.Range("G2").Value() + 7
$ws.Range("J2").Value = $ws.Range("J2").Value() + 1
$ws.Range("N2").Value = $ws.Range("N2").Value() + 1
$ws.Range("O2").Value = $ws.Range("O2").Value() + 1
$ws.Range("P2").Value = $ws.Range("P2").Value() + 1

$ws.Range("B3").Value = $ws.Range("B3").Value() + 1
$ws.Range("C3").Value = $ws.Range("C3").Value() + 15
$ws.Range("E3").Value = $ws.Range("E3").Value() + 2
$ws.Range("F3").Value = $ws.Range("F3").Value() + 8
$ws.Range("G3").Value = $ws.Range("G3").Value() + 1
$ws.Range("H3").Value = $ws.Range("H3").Value() + 2
$ws.Range("I3").Value = $ws.Range("I3").Value() + 3
$ws.Range("J3").Value = $ws.Range("J3").Value() + 3
$ws.Range("L3").Value = $ws.Range("L3").Value() + 35
$ws.Range("M3").Value = $ws.Range("M3").Value() + 26
$ws.Range("Q3").Value = $ws.Range("Q3").Value() + 56

# ---------------------------------------------------------------------
# DEF sheet - bump Home (row 2) and Road (row 3) season totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")

$ws.Range("C2").Value = $ws.Range("C2").Value() + 14
$ws.Range("D2").Value = $ws.Range("D2").Value() + 1
$ws.Range("F2").Value = $ws.Range("F2").Value() + 3
$ws.Range("G2").Value = $ws.Range("G2").Value() + 7
$ws.Range("J2").Value = $ws.Range("J2").Value() + 5
$ws.Range("N2").Value = $ws.Range("N2").Value() + 5
$ws.Range("O2").Value = $ws.Range("O2").Value() + 3
$ws.Range("P2").Value = $ws.Range("P2").Value() + 3

$ws.Range("B3").Value = $ws.Range("B3").Value() + 1
$ws.Range("C3").Value = $ws.Range("C3").Value() + 11
$ws.Range("E3").Value = $ws.Range("E3").Value() + 2
$ws.Range("F3").Value = $ws.Range("F3").Value() + 8
$ws.Range("G3").Value = $ws.Range("G3").Value() + 3
$ws.Range("H3").Value = $ws.Range("H3").Value() + 2
$ws.Range("I3").Value = $ws.Range("I3").Value() + 2
$ws.Range("J3").Value = $ws.Range("J3").Value() + 4
$ws.Range("L3").Value = $ws.Range("L3").Value() + 29
$ws.Range("M3").Value = $ws.Range("M3").Value() + 19
$ws.Range("Q3").Value = $ws.Range("Q3").Value() + 69

# ---------------------------------------------------------------------
# TURNS sheet - fix Road (row 3) turnover totals (tiebreak data fix)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")

$ws.Range("B3").Value = $ws.Range("B3").Value() + 2
$ws.Range("C3").Value = $ws.Range("C3").Value() + 1
$ws.Range("E3").Value = $ws.Range("E3").Value() - 3

# ---------------------------------------------------------------------
# PEN sheet - bump penalty totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")

$ws.Range("D2").Value = $ws.Range("D2").Value() + 1
$ws.Range("B3").Value = $ws.Range("B3").Value() + 1
